$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row cell updates (Coin / Link / Price / Volume(1h)) pulled from the refreshed
# coinranking.com snapshot. Price values that look like plain decimals must be
# forced to text (NumberFormat "@") before assignment so Excel does not silently
# convert "27.06" -> 27.059999999999999 / drop the trailing zero in "163.20".
$updates = @(
    @{ Row = 2; D = "75.068.75"; E = "  +2.44%  " },
    @{ Row = 3; D = "2.814.61"; E = "  +8.53%  " },
    @{ Row = 4; E = "  +0.17%  " },
    @{ Row = 5; D = "188.94"; E = "  +3.35%  " },
    @{ Row = 6; D = "595.73"; E = "  +2.78%  " },
    @{ Row = 7; E = "  +0.10%  " },
    @{ Row = 8; E = "  +4.41%  " },
    @{ Row = 9; D = "0.194"; E = "  -0.16%  " },
    @{ Row = 10; D = "2.812.13"; E = "  +8.47%  " },
    @{ Row = 11; E = "  -0.57%  " },
    @{ Row = 12; E = "  +3.46%  " },
    @{ Row = 13; D = "4.83"; E = "  +2.24%  " },
    @{ Row = 14; D = "3.331.54"; E = "  +8.85%  " },
    @{ Row = 15; D = "74.959.25"; E = "  +2.46%  " },
    @{ Row = 16; E = "  +1.98%  " },
    @{ Row = 17; D = "27.06"; E = "  +5.02%  " },
    @{ Row = 18; D = "2.815.18"; E = "  +8.87%  " },
    @{ Row = 19; D = "9.09"; E = "  +3.84%  " },
    @{ Row = 20; D = "12.34"; E = "  +5.58%  " },
    @{ Row = 21; D = "376.28"; E = "  +2.29%  " },
    @{ Row = 22; E = "  +2.06%  " },
    @{ Row = 23; D = "4.10"; E = "  +1.22%  " },
    @{ Row = 24; D = "6.18"; E = "  -0.04%  " },
    @{ Row = 25; E = "  -0.31%  " },
    @{ Row = 26; D = "70.96"; E = "  +2.55%  " },
    @{ Row = 27; B = "WrappedeETH"; C = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"; D = "2.976.37"; E = "  +9.34%  " },
    @{ Row = 28; B = "NEARProtocol"; C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D = "4.19"; E = "  +2.50%  " },
    @{ Row = 29; D = "9.59"; E = "  +4.74%  " },
    @{ Row = 30; D = "0.0000104"; E = "  +12.82%  " },
    @{ Row = 31; D = "0.997"; E = "  -0.71%  " },
    @{ Row = 32; D = "516.37"; E = "  +3.46%  " },
    @{ Row = 33; E = "  +1.83%  " },
    @{ Row = 34; D = "7.82"; E = "  +3.01%  " },
    @{ Row = 35; E = "  +4.53%  " },
    @{ Row = 36; E = "  +0.05%  " },
    @{ Row = 37; D = "163.20"; E = "  +1.93%  " },
    @{ Row = 38; D = "20.10"; E = "  +5.96%  " },
    @{ Row = 39; E = "  -0.29%  " },
    @{ Row = 40; D = "19.33"; E = "  +0.58%  " },
    @{ Row = 41; D = "184.82"; E = "  +19.86%  " },
    @{ Row = 42; E = "  +0.05%  " },
    @{ Row = 43; E = "  +5.22%  " },
    @{ Row = 44; D = "0.341"; E = "  +5.69%  " },
    @{ Row = 45; E = "  +1.55%  " },
    @{ Row = 46; E = "  +5.50%  " },
    @{ Row = 47; D = "40.13"; E = "  +4.62%  " },
    @{ Row = 48; D = "2.35"; E = "  +1.49%  " },
    @{ Row = 49; E = "  +0.52%  " },
    @{ Row = 50; D = "0.569"; E = "  +9.43%  " },
    @{ Row = 51; D = "3.73"; E = "  +3.95%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Range("B" + $u.Row).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C" + $u.Row).Value = $u.C }
    if ($u.ContainsKey("D")) {
        $priceCell = $ws.Range("D" + $u.Row)
        if ($u.D -match "^[+-]?\d+(\.\d+)?([eE][+-]?\d+)?$") {
            $priceCell.NumberFormat = "@"
        }
        $priceCell.Value = $u.D
    }
    if ($u.ContainsKey("E")) { $ws.Range("E" + $u.Row).Value = $u.E }
}
